$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "36.884.67"
$ws.Cells.Item(2, 5).Value = "  -2.02%  "

$ws.Cells.Item(3, 4).Value = "1.996.75"
$ws.Cells.Item(3, 5).Value = "  -3.93%  "

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "0.995"
$c.ClearFormats()
$ws.Cells.Item(4, 5).Value = "  -0.65%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "224.01"
$c.ClearFormats()
$ws.Cells.Item(5, 5).Value = "  -3.62%  "

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "0.603"
$c.ClearFormats()
$ws.Cells.Item(6, 5).Value = "  -3.28%  "

$ws.Cells.Item(7, 5).Value = "  +0.01%  "

$ws.Cells.Item(8, 5).Value = "  -6.87%  "

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.374"
$c.ClearFormats()
$ws.Cells.Item(9, 5).Value = "  -3.58%  "

$ws.Cells.Item(10, 5).Value = "  -0.06%  "

$ws.Cells.Item(11, 5).Value = "  -4.82%  "

$ws.Cells.Item(12, 4).Value = "2.291.78"
$ws.Cells.Item(12, 5).Value = "  -3.88%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "14.00"
$c.ClearFormats()
$ws.Cells.Item(13, 5).Value = "  -5.56%  "

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "19.99"
$c.ClearFormats()
$ws.Cells.Item(14, 5).Value = "  -6.41%  "

$ws.Cells.Item(15, 5).Value = "  -4.81%  "

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "5.05"
$c.ClearFormats()
$ws.Cells.Item(16, 5).Value = "  -5.91%  "

$ws.Cells.Item(17, 4).Value = "1.982.57"
$ws.Cells.Item(17, 5).Value = "  -4.57%  "

$ws.Cells.Item(18, 4).Value = "36.816.28"
$ws.Cells.Item(18, 5).Value = "  -2.12%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "6.06"
$c.ClearFormats()
$ws.Cells.Item(19, 5).Value = "  -1.60%  "

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "68.34"
$c.ClearFormats()
$ws.Cells.Item(20, 5).Value = "  -2.42%  "

$ws.Cells.Item(21, 4).Value = "0.0₃0809"
$ws.Cells.Item(21, 5).Value = "  -2.36%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "221.82"
$c.ClearFormats()
$ws.Cells.Item(22, 5).Value = "  -2.60%  "

$ws.Cells.Item(23, 5).Value = "  -0.07%  "

$ws.Cells.Item(24, 5).Value = "  +1.66%  "

$ws.Cells.Item(25, 5).Value = "  -8.45%  "

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "165.03"
$c.ClearFormats()
$ws.Cells.Item(26, 5).Value = "  -2.88%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "9.06"
$c.ClearFormats()
$ws.Cells.Item(27, 5).Value = "  -8.37%  "

$ws.Cells.Item(28, 5).Value = "  -3.93%  "

$ws.Cells.Item(29, 5).Value = "  -2.14%  "

$ws.Cells.Item(30, 5).Value = "  -7.65%  "

$ws.Cells.Item(31, 5).Value = "  -4.76%  "

$ws.Cells.Item(32, 5).Value = "  -2.51%  "

$ws.Cells.Item(33, 5).Value = "  -3.54%  "

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "4.37"
$c.ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -5.89%  "

$ws.Cells.Item(35, 5).Value = "  -8.72%  "

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "1.84"
$c.ClearFormats()
$ws.Cells.Item(36, 5).Value = "  +0.90%  "

$ws.Cells.Item(37, 5).Value = "  -0.02%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "3.12"
$c.ClearFormats()
$ws.Cells.Item(38, 5).Value = "  -5.81%  "

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "5.22"
$c.ClearFormats()
$ws.Cells.Item(39, 5).Value = "  -1.92%  "

$ws.Cells.Item(40, 4).Value = "1.461.78"
$ws.Cells.Item(40, 5).Value = "  -1.39%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.0214"
$c.ClearFormats()
$ws.Cells.Item(41, 5).Value = "  -6.06%  "

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "94.06"
$c.ClearFormats()
$ws.Cells.Item(42, 5).Value = "  -5.04%  "

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "0.0911"
$c.ClearFormats()
$ws.Cells.Item(43, 5).Value = "  -5.29%  "

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "16.12"
$c.ClearFormats()
$ws.Cells.Item(44, 5).Value = "  -4.15%  "

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "2.75"
$c.ClearFormats()
$ws.Cells.Item(45, 5).Value = "  -5.35%  "

$ws.Cells.Item(46, 5).Value = "  -7.60%  "

$ws.Cells.Item(47, 2).Value = "ARBITRUM"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Cells.Item(47, 5).Value = "  -3.93%  "

$ws.Cells.Item(48, 2).Value = "FraxShare"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "7.10"
$c.ClearFormats()
$ws.Cells.Item(48, 5).Value = "  -2.23%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "2.90"
$c.ClearFormats()
$ws.Cells.Item(49, 5).Value = "  -2.66%  "

$ws.Cells.Item(50, 4).Value = "2.180.42"
$ws.Cells.Item(50, 5).Value = "  -3.85%  "

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "43.90"
$c.ClearFormats()
$ws.Cells.Item(51, 5).Value = "  -4.45%  "
